# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" summary row at the top of the "总计" sheet and
#    renumber the existing index column.
# 2. Insert a brand-new "2022-Q4" worksheet (built from the "2022-Q3"
#    worksheet so it inherits the same formatting) with the Q4 fund-holding
#    detail data, positioned right after "总计" and before "2022-Q3".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" (summary) sheet - insert new row for 2022-Q4
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push all existing data rows down by one.
$summary.Rows.Item(2).Insert()

# Match formatting of the row below (style carries the bordered/centered
# index-column look + plain data cells) before writing the new values.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.18

# Renumber the index column (A) for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

# ---------------------------------------------------------------------
# Step 2: new "2022-Q4" worksheet with per-fund holding detail
# ---------------------------------------------------------------------
# "2022-Q3" (currently Worksheets.Item(2)) is used as a template so the
# new sheet inherits identical header/index-column styling; it is copied
# to a position right before itself, i.e. right after "总计".
$template = $wb.Worksheets.Item(2)
$template.Copy($template, $null)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template only has 7 rows (1 header + 6 data); Q4 needs 8 (1 header +
# 7 data), so add one more row, matching the formatting already in use.
$q4.Range("A7:H7").Copy()
$q4.Range("A8:H8").PasteSpecial(-4122)

# Columns B-G hold fund codes / names / size / position figures that must
# stay text (e.g. fund code "001009" would otherwise be auto-coerced to
# the number 1009, losing the leading zeros) - force text before writing.
$q4.Range("B2:G8").NumberFormat = "@"

# Row 2
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "001009"
$q4.Range("C2").Value = "上投摩根安全战略股票A"
$q4.Range("D2").Value = "3.32"
$q4.Range("E2").Value = "91.31"
$q4.Range("F2").Value = "2.83"
$q4.Range("G2").Value = "0.0940"
$q4.Range("H2").Value = 8

# Row 3
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "001484"
$q4.Range("C3").Value = "天弘新价值灵活配置混合A"
$q4.Range("D3").Value = "1.26"
$q4.Range("E3").Value = "94.05"
$q4.Range("F3").Value = "2.54"
$q4.Range("G3").Value = "0.0320"
$q4.Range("H3").Value = 9

# Row 4
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "006700"
$q4.Range("C4").Value = "红土创新稳健混合A"
$q4.Range("D4").Value = "0.51"
$q4.Range("E4").Value = "23.63"
$q4.Range("F4").Value = "4.71"
$q4.Range("G4").Value = "0.0240"
$q4.Range("H4").Value = 2

# Row 5
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "006701"
$q4.Range("C5").Value = "红土创新稳健混合C"
$q4.Range("D5").Value = "0.44"
$q4.Range("E5").Value = "23.63"
$q4.Range("F5").Value = "4.71"
$q4.Range("G5").Value = "0.0207"
$q4.Range("H5").Value = 2

# Row 6
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "011816"
$q4.Range("C6").Value = "融通多元收益一年持有期混合"
$q4.Range("D6").Value = "0.78"
$q4.Range("E6").Value = "22.61"
$q4.Range("F6").Value = "1.07"
$q4.Range("G6").Value = "0.0083"
$q4.Range("H6").Value = 4

# Row 7
$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "015172"
$q4.Range("C7").Value = "上投摩根安全战略股票C"
$q4.Range("D7").Value = "0.07"
$q4.Range("E7").Value = "91.31"
$q4.Range("F7").Value = "2.83"
$q4.Range("G7").Value = "0.0020"
$q4.Range("H7").Value = 8

# Row 8
$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "016246"
$q4.Range("C8").Value = "天弘新价值灵活配置混合C"
$q4.Range("D8").Value = "0.03"
$q4.Range("E8").Value = "94.05"
$q4.Range("F8").Value = "2.54"
$q4.Range("G8").Value = "0.0008"
$q4.Range("H8").Value = 9
